$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    current 2nd sheet, which is "2022-Q3"). All the other quarter sheets
#    shift one position to the right, which also happens automatically.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$newSheet.Name = "2022-Q4"

$totalSheet = $wb.Worksheets.Item(1)

# Fill in the header row (row 1) of the new sheet, matching the layout used
# by the other quarterly sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2
    $newSheet.Cells.Item(1, $col).Value = $headers[$i]
}
# Copy the header cell style (bold + border, centered) from the existing
# "总计" sheet header onto the new header row.
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Fund-level detail rows for 2022-Q4.
$q4Data = @(
    @(0, "009341", "易方达均衡成长股票",           "59.89", "88.68", "2.79", "1.6709", 7),
    @(1, "013554", "信澳远见价值混合A",             "0.86",  "67.74", "3.83", "0.0329", 5),
    @(2, "011471", "鹏华致远成长混合A",             "1.67",  "60.84", "1.79", "0.0299", 9),
    @(3, "000059", "国联安中证医药100指数A",        "1.93",  "92.58", "1.11", "0.0214", 3),
    @(4, "013555", "信澳远见价值混合C",             "0.46",  "67.74", "3.83", "0.0176", 5),
    @(5, "006569", "国联安中证医药100指数C",        "0.97",  "92.58", "1.11", "0.0108", 3),
    @(6, "011433", "中加聚优一年定期开放混合A",     "0.87",  "24.92", "1.24", "0.0108", 5),
    @(7, "005281", "中科沃土转型升级灵活配置混合",  "0.10",  "60.18", "3.88", "0.0039", 2),
    @(8, "011472", "鹏华致远成长混合C",             "0.06",  "60.84", "1.79", "0.0011", 9),
    @(9, "011434", "中加聚优一年定期开放混合C",     "0.03",  "24.92", "1.24", "0.0004", 5)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = $i + 2
    $rowData = $q4Data[$i]

    $newSheet.Cells.Item($row, 1).Value = $rowData[0]
    # Force text storage (not auto-numeric) for code/name/percent-like
    # columns by prefixing with an apostrophe, then strip the formatting
    # that Excel applies automatically when doing so.
    $newSheet.Cells.Item($row, 2).Value = "'" + $rowData[1]
    $newSheet.Cells.Item($row, 3).Value = $rowData[2]
    $newSheet.Cells.Item($row, 4).Value = "'" + $rowData[3]
    $newSheet.Cells.Item($row, 5).Value = "'" + $rowData[4]
    $newSheet.Cells.Item($row, 6).Value = "'" + $rowData[5]
    $newSheet.Cells.Item($row, 7).Value = "'" + $rowData[6]
    $newSheet.Cells.Item($row, 8).Value = $rowData[7]

    $newSheet.Range($newSheet.Cells.Item($row, 2), $newSheet.Cells.Item($row, 7)).ClearFormats()
}

# Copy the index-column style (bold border, centered) used in column A of
# the "总计" sheet onto column A of the new sheet's data rows.
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Prepend a 2022-Q4 summary row to the "总计" (total) sheet, shifting the
#    existing quarter rows down by one and renumbering the index column.
# ---------------------------------------------------------------------------
$summaryData = @(
    @(0, "2022-Q4", 10, 1.8),
    @(1, "2022-Q3", 1, 0),
    @(2, "2022-Q2", 2, 1.67),
    @(3, "2022-Q1", 3, 1.1),
    @(4, "2021-Q4", 2, 3.55),
    @(5, "2021-Q3", 7, 5.95),
    @(6, "2021-Q2", 2, 3.02),
    @(7, "2021-Q1", 11, 8.380000000000001)
)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $row = $i + 2
    $rowData = $summaryData[$i]
    $totalSheet.Cells.Item($row, 1).Value = $rowData[0]
    $totalSheet.Cells.Item($row, 2).Value = $rowData[1]
    $totalSheet.Cells.Item($row, 3).Value = $rowData[2]
    $totalSheet.Cells.Item($row, 4).Value = $rowData[3]
}

# The very last row (2021-Q1) used to be the bottom row already carrying the
# styled index-column format; after rewriting it directly via .Value it
# needs that format (re-)applied, same as the rest of column A.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)
